$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename Sheet1 to "Add a Loanbook data"
$ws.Name = "Add a Loanbook data"

# Update font color on the used range (theme color -> explicit black RGB)
$ws.Range("A1:G6").Font.Color = 0

# Update row heights (Excel re-measured heights after the font change)
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 186.75
$ws.Rows.Item(3).RowHeight = 20.25
$ws.Rows.Item(4).RowHeight = 20.25
$ws.Rows.Item(5).RowHeight = 20.25
$ws.Rows.Item(6).RowHeight = 21
